$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 101 (G101=19884)
$ws.Range("H101").Value = 750
$ws.Range("J101").Value = 1000
$ws.Range("L101").Value = 3000
$ws.Range("N101").Value = -6244
# Row 132 (G132=44049)
$ws.Range("H132").Value = 1236.5555
$ws.Range("I132").Value = 1236.5555
$ws.Range("K132").Value = 3709.6665
$ws.Range("M132").Value = -1179.6665
# Row 138 (G138=44169)
$ws.Range("H138").Value = 2558.6086
$ws.Range("I138").Value = 1382.8334
$ws.Range("J138").Value = 2973.5881
$ws.Range("K138").Value = 4148.5002
$ws.Range("L138").Value = 8920.764299999999
$ws.Range("M138").Value = 991.4997999999996
$ws.Range("N138").Value = -19200.7643

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 897.6429000000001
$ws.Range("I32").Value = 455.66666
$ws.Range("K32").Value = 455.66666
$ws.Range("M32").Value = -168.66666
# Row 45 (G45=27714)
$ws.Range("H45").Value = 2670.7144
$ws.Range("I45").Value = 2539
$ws.Range("K45").Value = 2539
$ws.Range("M45").Value = -2162
# Row 102 (G102=19945)
$ws.Range("H102").Value = 52501504
$ws.Range("I102").Value = 5001505
$ws.Range("K102").Value = 5001505
$ws.Range("M102").Value = -4999883
# Row 121 (G121=26285)
$ws.Range("H121").Value = 41500
$ws.Range("J121").Value = 41500
$ws.Range("L121").Value = 41500
$ws.Range("N121").Value = -44994
# Row 122 (G122=36168)
$ws.Range("H122").Value = 7365
$ws.Range("I122").Value = 8725
$ws.Range("J122").Value = 1925
$ws.Range("K122").Value = 26175
$ws.Range("L122").Value = 5775
$ws.Range("M122").Value = -23725
$ws.Range("N122").Value = -10675

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (G20=14149)
$ws.Range("H20").Value = 987.25
$ws.Range("J20").Value = 1350
$ws.Range("L20").Value = 1350
$ws.Range("N20").Value = -1844
# Row 76 (G76=10630)
$ws.Range("H76").Value = 18749.75
$ws.Range("J76").Value = 18749.75
$ws.Range("L76").Value = 18749.75
$ws.Range("N76").Value = -19379.75
# Row 79 (G79=10630)
$ws.Range("H79").Value = 18749.75
$ws.Range("J79").Value = 18749.75
$ws.Range("L79").Value = 18749.75
$ws.Range("N79").Value = -20933.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (G7=5361)
$ws.Range("H7").Value = 19.142857
$ws.Range("I7").Value = 19.142857
$ws.Range("K7").Value = 19.142857
$ws.Range("M7").Value = 93.85714300000001
# Row 31 (G31=44023)
$ws.Range("H31").Value = 1498.8572
$ws.Range("I31").Value = 1248.6666
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 1248.6666
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -953.6666
$ws.Range("N31").Value = -3590
# Row 34 (G34=44023)
$ws.Range("H34").Value = 1498.8572
$ws.Range("I34").Value = 1248.6666
$ws.Range("K34").Value = 1248.6666
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = -1046.6666
$ws.Range("N34").Value = -3404
# Row 58 (G58=44021)
$ws.Range("H58").Value = 2560.3125
$ws.Range("I58").Value = 1864.3334
$ws.Range("J58").Value = 13000
$ws.Range("K58").Value = 1864.3334
$ws.Range("L58").Value = 13000
$ws.Range("M58").Value = -1661.3334
$ws.Range("N58").Value = -13406
# Row 82 (G82=10799)
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
# Row 85 (G85=10799)
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
# Row 136 (G136=44021)
$ws.Range("H136").Value = 2560.3125
$ws.Range("I136").Value = 1864.3334
$ws.Range("J136").Value = 13000
$ws.Range("K136").Value = 5593.0002
$ws.Range("L136").Value = 39000
$ws.Range("M136").Value = -3043.0002
$ws.Range("N136").Value = -44100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 9 (G9=4681)
$ws.Range("H9").Value = 5940.8
$ws.Range("I9").Value = 4300
$ws.Range("K9").Value = 12900
$ws.Range("M9").Value = -12676
# Row 46 (G46=4701)
$ws.Range("H46").Value = 5075
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
# Row 111 (G111=27856)
$ws.Range("H111").Value = 1027
$ws.Range("I111").Value = 1027
$ws.Range("K111").Value = 3081
$ws.Range("M111").Value = -14
# Row 129 (G129=36054)
$ws.Range("H129").Value = 1398.2
$ws.Range("J129").Value = 1644.3334
$ws.Range("L129").Value = 4933.0002
$ws.Range("N129").Value = -14933.0002
# Row 131 (G131=36060)
$ws.Range("H131").Value = 2825
$ws.Range("I131").Value = 2535.8
$ws.Range("K131").Value = 7607.400000000001
$ws.Range("M131").Value = -2567.400000000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 92 (G92=18094)
$ws.Range("H92").Value = 54500
$ws.Range("J92").Value = 54500
$ws.Range("L92").Value = 54500
$ws.Range("N92").Value = -58244
# Row 97 (G97=19940)
$ws.Range("H97").Value = 1029.6
$ws.Range("J97").Value = 1677.75
$ws.Range("L97").Value = 1677.75
$ws.Range("N97").Value = -2669.75
# Row 102 (G102=36169)
$ws.Range("H102").Value = 892.2727
$ws.Range("I102").Value = 898.3333
$ws.Range("J102").Value = 865
$ws.Range("K102").Value = 898.3333
$ws.Range("L102").Value = 865
$ws.Range("M102").Value = 723.6667
$ws.Range("N102").Value = -4109
# Row 122 (G122=36182)
$ws.Range("H122").Value = 2821.5557
$ws.Range("I122").Value = 1879.4
$ws.Range("K122").Value = 5638.200000000001
$ws.Range("M122").Value = -3188.200000000001
# Row 126 (G126=36184)
$ws.Range("H126").Value = 1985
$ws.Range("I126").Value = 1985
$ws.Range("K126").Value = 5955
$ws.Range("M126").Value = -3485
# Row 132 (G132=44008)
$ws.Range("H132").Value = 1561.2858
$ws.Range("J132").Value = 2500
$ws.Range("L132").Value = 7500
$ws.Range("N132").Value = -12560

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 5 (G5=3790)
$ws.Range("H5").Value = 20000
$ws.Range("J5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("N5").Value = -20226
# Row 40 (G40=36248)
$ws.Range("H40").Value = 1672667.4
$ws.Range("I40").Value = 3003
$ws.Range("K40").Value = 3003
$ws.Range("M40").Value = -2867
# Row 61 (G61=27740)
$ws.Range("H61").Value = 2849.6
$ws.Range("I61").Value = 2849.6
$ws.Range("K61").Value = 2849.6
$ws.Range("M61").Value = -2647.6
# Row 68 (G68=12563)
$ws.Range("H68").Value = 2000
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
# Row 71 (G71=12563)
$ws.Range("H71").Value = 2000
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
# Row 100 (G100=19995)
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
# Row 113 (G113=27740)
$ws.Range("H113").Value = 2849.6
$ws.Range("I113").Value = 2849.6
$ws.Range("K113").Value = 2849.6
$ws.Range("M113").Value = -679.5999999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 105 (G105=18710)
$ws.Range("H105").Value = 45750
$ws.Range("J105").Value = 45750
$ws.Range("L105").Value = 45750
$ws.Range("N105").Value = -52738
# Row 113 (G113=27752)
$ws.Range("H113").Value = 8219
$ws.Range("I113").Value = 422.16666
$ws.Range("J113").Value = 55000
$ws.Range("K113").Value = 1266.49998
$ws.Range("L113").Value = 165000
$ws.Range("M113").Value = 903.5000199999999
$ws.Range("N113").Value = -169340
